$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "up-to-mid cap" values for column C (rows 2,4-20; row 1 and row 3 keep "80").
$values = @{
    2  = "48.259"
    4  = "22.171"
    5  = "224.341"
    6  = "139.261"
    7  = "123.703"
    8  = "200.571"
    9  = "112.849"
    10 = "151.427"
    11 = "5.69"
    12 = "0.617"
    13 = "319.289"
    14 = "253.039"
    15 = "250.311"
    16 = "410.517"
    17 = "175.847"
    18 = "125.506"
    19 = "149.405"
    20 = "147.022"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("C$row")
    # Force text storage so the numeric-looking string isn't auto-converted
    # to a real number, matching the existing column C text cells.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$row]
    # Drop the temporary formatting again so the cell keeps the sheet's
    # default (General) style, same as all the other cells.
    $cell.ClearFormats()
}
